$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '61.147.78'
$c.Style = 'Normal'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  -4.19%  '
$c.Style = 'Normal'

$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '2.976.63'
$c.Style = 'Normal'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  -3.59%  '
$c.Style = 'Normal'

$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.Style = 'Normal'

$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '545.67'
$c.Style = 'Normal'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  +0.52%  '
$c.Style = 'Normal'

$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '132.35'
$c.Style = 'Normal'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -4.57%  '
$c.Style = 'Normal'

$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +0.06%  '
$c.Style = 'Normal'

$c = $ws.Range("D8")
$c.NumberFormat = '@'
$c.Value = '2.972.93'
$c.Style = 'Normal'
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  -3.54%  '
$c.Style = 'Normal'

$c = $ws.Range("D9")
$c.NumberFormat = '@'
$c.Value = '0.491'
$c.Style = 'Normal'
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  -0.94%  '
$c.Style = 'Normal'

$c = $ws.Range("D10")
$c.NumberFormat = '@'
$c.Value = '0.146'
$c.Style = 'Normal'
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  -6.02%  '
$c.Style = 'Normal'

$c = $ws.Range("D11")
$c.NumberFormat = '@'
$c.Value = '5.94'
$c.Style = 'Normal'
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  -8.78%  '
$c.Style = 'Normal'

$c = $ws.Range("D12")
$c.NumberFormat = '@'
$c.Value = '0.445'
$c.Style = 'Normal'
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  -2.38%  '
$c.Style = 'Normal'

$c = $ws.Range("B13")
$c.NumberFormat = '@'
$c.Value = 'ShibaInu'
$c.Style = 'Normal'
$c = $ws.Range("C13")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c.Style = 'Normal'
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '0.0000219'
$c.Style = 'Normal'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  -3.22%  '
$c.Style = 'Normal'

$c = $ws.Range("B14")
$c.NumberFormat = '@'
$c.Value = 'Avalanche'
$c.Style = 'Normal'
$c = $ws.Range("C14")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c.Style = 'Normal'
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '33.98'
$c.Style = 'Normal'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  -2.05%  '
$c.Style = 'Normal'

$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '3.468.49'
$c.Style = 'Normal'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  -3.35%  '
$c.Style = 'Normal'

$c = $ws.Range("B16")
$c.NumberFormat = '@'
$c.Value = 'TRON'
$c.Style = 'Normal'
$c = $ws.Range("C16")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c.Style = 'Normal'
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '0.110'
$c.Style = 'Normal'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -2.66%  '
$c.Style = 'Normal'

$c = $ws.Range("B17")
$c.NumberFormat = '@'
$c.Value = 'WrappedBTC'
$c.Style = 'Normal'
$c = $ws.Range("C17")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$c.Style = 'Normal'
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '61.279.08'
$c.Style = 'Normal'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  -4.15%  '
$c.Style = 'Normal'

$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '2.982.60'
$c.Style = 'Normal'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  -3.51%  '
$c.Style = 'Normal'

$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '6.60'
$c.Style = 'Normal'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -0.57%  '
$c.Style = 'Normal'

$c = $ws.Range("D20")
$c.NumberFormat = '@'
$c.Value = '467.79'
$c.Style = 'Normal'
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -2.41%  '
$c.Style = 'Normal'

$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '13.12'
$c.Style = 'Normal'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  -2.20%  '
$c.Style = 'Normal'

$c = $ws.Range("D22")
$c.NumberFormat = '@'
$c.Value = '0.668'
$c.Style = 'Normal'
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  -4.32%  '
$c.Style = 'Normal'

$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '6.98'
$c.Style = 'Normal'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  -1.47%  '
$c.Style = 'Normal'

$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '79.80'
$c.Style = 'Normal'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  +1.24%  '
$c.Style = 'Normal'

$c = $ws.Range("D25")
$c.NumberFormat = '@'
$c.Value = '11.98'
$c.Style = 'Normal'
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  -2.78%  '
$c.Style = 'Normal'

$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.Style = 'Normal'

$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '2.70'
$c.Style = 'Normal'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  -0.85%  '
$c.Style = 'Normal'

$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '7.67'
$c.Style = 'Normal'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  -4.80%  '
$c.Style = 'Normal'

$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  -0.14%  '
$c.Style = 'Normal'

$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '1.89'
$c.Style = 'Normal'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -0.62%  '
$c.Style = 'Normal'

$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '25.39'
$c.Style = 'Normal'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  -3.44%  '
$c.Style = 'Normal'

$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '1.12'
$c.Style = 'Normal'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  -3.41%  '
$c.Style = 'Normal'

$c = $ws.Range("B33")
$c.NumberFormat = '@'
$c.Value = 'Stacks'
$c.Style = 'Normal'
$c = $ws.Range("C33")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c.Style = 'Normal'
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '2.29'
$c.Style = 'Normal'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  -2.57%  '
$c.Style = 'Normal'

$c = $ws.Range("B34")
$c.NumberFormat = '@'
$c.Value = 'NEARProtocol'
$c.Style = 'Normal'
$c = $ws.Range("C34")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$c.Style = 'Normal'
$c = $ws.Range("D34")
$c.NumberFormat = '@'
$c.Value = '5.44'
$c.Style = 'Normal'
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +1.82%  '
$c.Style = 'Normal'

$c = $ws.Range("D35")
$c.NumberFormat = '@'
$c.Value = '55.13'
$c.Style = 'Normal'
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  -3.34%  '
$c.Style = 'Normal'

$c = $ws.Range("D36")
$c.NumberFormat = '@'
$c.Value = '5.84'
$c.Style = 'Normal'
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  -2.69%  '
$c.Style = 'Normal'

$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '445.44'
$c.Style = 'Normal'
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -9.55%  '
$c.Style = 'Normal'

$c = $ws.Range("D38")
$c.NumberFormat = '@'
$c.Value = '3.153.43'
$c.Style = 'Normal'
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  -2.85%  '
$c.Style = 'Normal'

$c = $ws.Range("D39")
$c.NumberFormat = '@'
$c.Value = '0.0789'
$c.Style = 'Normal'
$c = $ws.Range("E39")
$c.NumberFormat = '@'
$c.Value = '  -1.10%  '
$c.Style = 'Normal'

$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '0.0378'
$c.Style = 'Normal'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -6.28%  '
$c.Style = 'Normal'

$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  -2.48%  '
$c.Style = 'Normal'

$c = $ws.Range("D42")
$c.NumberFormat = '@'
$c.Value = '8.08'
$c.Style = 'Normal'
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  -0.11%  '
$c.Style = 'Normal'

$c = $ws.Range("D43")
$c.NumberFormat = '@'
$c.Value = '2.39'
$c.Style = 'Normal'
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  -10.53%  '
$c.Style = 'Normal'

$c = $ws.Range("B44")
$c.NumberFormat = '@'
$c.Value = 'USDe'
$c.Style = 'Normal'
$c = $ws.Range("C44")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$c.Style = 'Normal'
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c.Style = 'Normal'

$c = $ws.Range("B45")
$c.NumberFormat = '@'
$c.Value = 'InjectiveProtocol'
$c.Style = 'Normal'
$c = $ws.Range("C45")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c.Style = 'Normal'
$c = $ws.Range("D45")
$c.NumberFormat = '@'
$c.Value = '25.76'
$c.Style = 'Normal'
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  +2.40%  '
$c.Style = 'Normal'

$c = $ws.Range("D46")
$c.NumberFormat = '@'
$c.Value = '0.241'
$c.Style = 'Normal'
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  -4.83%  '
$c.Style = 'Normal'

$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '0.107'
$c.Style = 'Normal'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -1.45%  '
$c.Style = 'Normal'

$c = $ws.Range("B48")
$c.NumberFormat = '@'
$c.Value = 'Monero'
$c.Style = 'Normal'
$c = $ws.Range("C48")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c.Style = 'Normal'
$c = $ws.Range("D48")
$c.NumberFormat = '@'
$c.Value = '116.97'
$c.Style = 'Normal'
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  -5.67%  '
$c.Style = 'Normal'

$c = $ws.Range("B49")
$c.NumberFormat = '@'
$c.Value = 'Fetch.AI'
$c.Style = 'Normal'
$c = $ws.Range("C49")
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$c.Style = 'Normal'
$c = $ws.Range("D49")
$c.NumberFormat = '@'
$c.Value = '1.95'
$c.Style = 'Normal'
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  -4.23%  '
$c.Style = 'Normal'

$c = $ws.Range("D50")
$c.NumberFormat = '@'
$c.Value = '1.29'
$c.Style = 'Normal'
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  +7.27%  '
$c.Style = 'Normal'

$c = $ws.Range("D51")
$c.NumberFormat = '@'
$c.Value = '0.0₃0484'
$c.Style = 'Normal'
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  -8.82%  '
$c.Style = 'Normal'
